$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New row 134: Harris poll
$ws.Cells.Item(134, 1).Value = 40
$ws.Cells.Item(134, 2).Value = 2021
$ws.Cells.Item(134, 3).Value = 14
$ws.Cells.Item(134, 4).Value = 12
$ws.Cells.Item(134, 5).Value = 5
$ws.Cells.Item(134, 6).Value = "harris"
$ws.Cells.Item(134, 7).Value = "online"
$ws.Cells.Item(134, 8).Value = "included"
$ws.Cells.Item(134, 9).Value = 2129
$ws.Cells.Item(134, 10).Value = 1
$ws.Cells.Item(134, 11).Value = 1
$ws.Cells.Item(134, 12).Value = 11
$ws.Cells.Item(134, 13).Value = 2
$ws.Cells.Item(134, 14).Value = 1
$ws.Cells.Item(134, 15).Value = 7
$ws.Cells.Item(134, 16).Value = 5
$ws.Cells.Item(134, 17).Value = 23
$ws.Cells.Item(134, 18).Value = 14
$ws.Cells.Item(134, 21).Value = "T_0.5"
$ws.Cells.Item(134, 22).Value = 2
$ws.Cells.Item(134, 23).Value = 18
$ws.Cells.Item(134, 24).Value = 14
$ws.Cells.Item(134, 25).Value = "T_0.5"
$ws.Cells.Item(134, 27).Value = 1

# New row 135: Ifop poll
$ws.Cells.Item(135, 1).Value = 41
$ws.Cells.Item(135, 2).Value = 2021
$ws.Cells.Item(135, 3).Value = 14
$ws.Cells.Item(135, 4).Value = 12
$ws.Cells.Item(135, 5).Value = 5
$ws.Cells.Item(135, 6).Value = "ifop"
$ws.Cells.Item(135, 7).Value = "online"
$ws.Cells.Item(135, 8).Value = "included"
$ws.Cells.Item(135, 9).Value = 1341
$ws.Cells.Item(135, 10).Value = 0.5
$ws.Cells.Item(135, 11).Value = 0.5
$ws.Cells.Item(135, 12).Value = 9
$ws.Cells.Item(135, 13).Value = 2.5
$ws.Cells.Item(135, 14).Value = 1.5
$ws.Cells.Item(135, 15).Value = 6
$ws.Cells.Item(135, 16).Value = 5
$ws.Cells.Item(135, 17).Value = 25
$ws.Cells.Item(135, 18).Value = 17
$ws.Cells.Item(135, 21).Value = 0.5
$ws.Cells.Item(135, 22).Value = 2.5
$ws.Cells.Item(135, 23).Value = 17
$ws.Cells.Item(135, 24).Value = 13

$ws.Range("Y134").Font.Color = 0

$ws.Range("AB135").Select()
$excel.ActiveWindow.ScrollRow = 126
$excel.ActiveWindow.ScrollColumn = 5

